# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Price column (D) and Volume(1h) column (E) are plain text in this sheet
# (e.g. "63.021.77" uses '.' as both thousands- and would-be decimal
# separator, and E holds space-padded "  +0.12%  " strings), so values
# that would otherwise be auto-parsed as a number are written with a
# leading single-quote to force Excel to keep them as text, matching the
# source data's string type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.021.77'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.547.43'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''583.48'
$ws.Range('E5').Value = '  +2.13%  '
$ws.Range('D6').Value = '''146.84'
$ws.Range('E6').Value = '  -2.53%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.584'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('D10').Value = '''5.56'
$ws.Range('E10').Value = '  -3.57%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').Value = '''0.354'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').Value = '''27.49'
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('D14').Value = '3.006.31'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').Value = '62.922.69'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '2.543.71'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '''11.34'
$ws.Range('E18').Value = '  -3.10%  '
$ws.Range('D19').Value = '''337.89'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').Value = '''4.33'
$ws.Range('E20').Value = '  -0.98%  '
$ws.Range('D21').Value = '''6.76'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = '''65.65'
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D24').Value = '2.678.63'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('B28').Value = 'SuiNetwork'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D28').Value = '''1.48'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').Value = '''8.35'
$ws.Range('E29').Value = '  -3.65%  '
$ws.Range('D30').Value = '''7.69'
$ws.Range('E30').Value = '  +6.81%  '
$ws.Range('E31').Value = '  +4.79%  '
$ws.Range('D32').Value = '0.0₃0814'
$ws.Range('E32').Value = '  -2.36%  '
$ws.Range('D33').Value = '''178.02'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''1.55'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').Value = '''417.28'
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('D36').Value = '''0.400'
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('D37').Value = '''19.12'
$ws.Range('E37').Value = '  -0.48%  '
$ws.Range('D39').Value = '''4.36'
$ws.Range('E39').Value = '  -2.43%  '
$ws.Range('E40').Value = '  -2.28%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = '''39.79'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '''150.68'
$ws.Range('E43').Value = '  -2.37%  '
$ws.Range('D44').Value = '''3.77'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('D45').Value = '''20.75'
$ws.Range('E45').Value = '  -2.04%  '
$ws.Range('D46').Value = '''0.0540'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('D47').Value = '''0.602'
$ws.Range('D48').Value = '''0.0969'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '''0.0238'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '''18.28'
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('D51').Value = '''1.71'
$ws.Range('E51').Value = '  -6.16%  '
